$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '47.285.49'
$ws.Range('E2').Value = '  +4.70%  '

Set-TextValue $ws 'D3' '2.640.51'
$ws.Range('E3').Value = '  +10.97%  '

$ws.Range('E4').Value = '  -0.06%  '

Set-TextValue $ws 'D5' '313.90'
$ws.Range('E5').Value = '  +7.32%  '

Set-TextValue $ws 'D6' '105.09'
$ws.Range('E6').Value = '  +11.97%  '

Set-TextValue $ws 'D7' '0.614'
$ws.Range('E7').Value = '  +10.88%  '

Set-TextValue $ws 'D8' '1.00'
$ws.Range('E8').Value = '  -0.03%  '

Set-TextValue $ws 'D9' '0.602'
$ws.Range('E9').Value = '  +21.53%  '

Set-TextValue $ws 'D10' '39.91'
$ws.Range('E10').Value = '  +18.05%  '

Set-TextValue $ws 'D11' '55.39'
$ws.Range('E11').Value = '  +3.68%  '

Set-TextValue $ws 'D12' '0.0858'
$ws.Range('E12').Value = '  +11.14%  '

Set-TextValue $ws 'D13' '8.47'
$ws.Range('E13').Value = '  +22.39%  '

Set-TextValue $ws 'D14' '3.045.74'
$ws.Range('E14').Value = '  +10.91%  '

$ws.Range('E15').Value = '  +3.61%  '

Set-TextValue $ws 'D16' '2.636.09'
$ws.Range('E16').Value = '  +10.75%  '

Set-TextValue $ws 'D17' '0.950'
$ws.Range('E17').Value = '  +16.19%  '

Set-TextValue $ws 'D18' '15.41'
$ws.Range('E18').Value = '  +11.19%  '

Set-TextValue $ws 'D19' '47.783.97'
$ws.Range('E19').Value = '  +5.92%  '

Set-TextValue $ws 'D20' '0.0000104'
$ws.Range('E20').Value = '  +12.18%  '

Set-TextValue $ws 'D21' '13.35'
$ws.Range('E21').Value = '  +8.68%  '

Set-TextValue $ws 'D22' '6.82'
$ws.Range('E22').Value = '  +12.35%  '

Set-TextValue $ws 'D23' '72.98'
$ws.Range('E23').Value = '  +10.59%  '

Set-TextValue $ws 'D24' '272.88'
$ws.Range('E24').Value = '  +14.97%  '

Set-TextValue $ws 'D25' '3.13'
$ws.Range('E25').Value = '  +14.22%  '

Set-TextValue $ws 'D26' '2.24'
$ws.Range('E26').Value = '  +20.22%  '

Set-TextValue $ws 'D27' '30.58'
$ws.Range('E27').Value = '  +46.36%  '

$ws.Range('E28').Value = '  -0.06%  '

$ws.Range('E29').Value = '  +1.27%  '

Set-TextValue $ws 'D30' '10.79'
$ws.Range('E30').Value = '  +13.79%  '

Set-TextValue $ws 'D31' '40.25'
$ws.Range('E31').Value = '  +8.01%  '

Set-TextValue $ws 'D32' '2.32'
$ws.Range('E32').Value = '  +4.91%  '

Set-TextValue $ws 'D33' '6.23'
$ws.Range('E33').Value = '  +16.43%  '

$ws.Range('E34').Value = '  -0.84%  '

$ws.Range('E35').Value = '  +17.03%  '

Set-TextValue $ws 'D36' '0.0858'
$ws.Range('E36').Value = '  +13.85%  '

$ws.Range('E37').Value = '  +6.42%  '

Set-TextValue $ws 'D38' '152.31'
$ws.Range('E38').Value = '  +3.54%  '

$ws.Range('E39').Value = '  +11.71%  '

Set-TextValue $ws 'D40' '0.126'
$ws.Range('E40').Value = '  +10.81%  '

$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D41' '23.65'
$ws.Range('E41').Value = '  +58.54%  '

$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws 'D42' '16.58'
$ws.Range('E42').Value = '  +13.54%  '

$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D43' '3.78'
$ws.Range('E43').Value = '  +20.42%  '

$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D44' '4.30'
$ws.Range('E44').Value = '  +16.36%  '

Set-TextValue $ws 'D45' '0.0335'
$ws.Range('E45').Value = '  +15.02%  '

Set-TextValue $ws 'D46' '2.179.03'
$ws.Range('E46').Value = '  +11.44%  '

Set-TextValue $ws 'D47' '95.84'
$ws.Range('E47').Value = '  +8.04%  '

Set-TextValue $ws 'D48' '0.997'
$ws.Range('E48').Value = '  -0.12%  '

Set-TextValue $ws 'D49' '10.06'
$ws.Range('E49').Value = '  +19.64%  '

Set-TextValue $ws 'D50' '114.79'
$ws.Range('E50').Value = '  +16.10%  '

$ws.Range('E51').Value = '  +6.92%  '
